# Add 7 new LHR Airports Ltd bond rows (rows 13-19) to the "Search Results"
# sheet, mirroring the existing Issuer/Coupon/Maturity/Issue Date/ISIN
# layout used by the other rows in the Master List of Bonds.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$issuer = "LHR Airports Ltd"

# Coupon, Maturity date, Issue date, ISIN
$data = @(
    @(6.375, "08/04/2028", "08/04/1998", "XS0089000516"),
    @(8.5,   "03/29/2021", "01/31/1996", "XS0063290711"),
    @(11.75, "03/31/2016", "02/28/1991", "XS0030487051"),
    @(5.75,  "12/10/2031", "12/10/2001", "XS0138797021"),
    @(5.75,  "12/10/2031", "12/10/2001", "XS0142079028"),
    @(5.75,  "11/27/2013", "11/27/2003", "XS0181263202"),
    @(5.125, "02/15/2023", "02/15/2006", "XS0243520722")
)

$firstRow = 13

$row = $firstRow
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $issuer
    $ws.Cells.Item($row, 3).Value = $item[0]
    $ws.Cells.Item($row, 6).Value = $item[3]
    $row = $row + 1
}

# Copy the date formatting (style) already used in column D/E from the last
# existing data row (12) down onto the new rows before writing the dates,
# so the new cells reuse the workbook's existing date style instead of
# creating a new one.
$lastRow = $firstRow + $data.Length - 1
$ws.Range("D12:E12").Copy()
$ws.Range("D" + $firstRow + ":E" + $lastRow).PasteSpecial(-4122)

$row = $firstRow
foreach ($item in $data) {
    $ws.Cells.Item($row, 4).Value2 = $item[1]
    $ws.Cells.Item($row, 5).Value2 = $item[2]
    $row = $row + 1
}

# Match the author's final selection/scroll position.
$null = $ws.Range("AY11").Select()
